$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 7).Value2 = 38.38922700000001
$ws.Cells.Item(2, 8).Value2 = 76.77845400000001
$ws.Cells.Item(2, 9).Value2 = 0.4452295445267456
$ws.Cells.Item(2, 10).Value2 = 0.3751312191747254
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 13).Value2 = 0.431063
$ws.Cells.Item(2, 14).Value2 = 1.293189
$ws.Cells.Item(2, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(2, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(2, 17).Value2 = 16.548175358301
$ws.Cells.Item(2, 18).Value2 = 99.28905214980601
$ws.Cells.Item(2, 19).Value2 = 0.004756564822150545
$ws.Cells.Item(2, 20).Value2 = 0.004010189888701697

$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 7).Value2 = 38.38922700000001
$ws.Cells.Item(3, 8).Value2 = 76.77845400000001
$ws.Cells.Item(3, 9).Value2 = 0.4452295445267456
$ws.Cells.Item(3, 10).Value2 = 0.3751312191747254
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 18.65163266666667
$ws.Cells.Item(3, 14).Value2 = 55.954898
$ws.Cells.Item(3, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(3, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(3, 17).Value2 = 716.0217603612822
$ws.Cells.Item(3, 18).Value2 = 4296.130562167693
$ws.Cells.Item(3, 19).Value2 = 0.2058114470922827
$ws.Cells.Item(3, 20).Value2 = 0.173516605989484

$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 7).Value2 = 38.38922700000001
$ws.Cells.Item(4, 8).Value2 = 76.77845400000001
$ws.Cells.Item(4, 9).Value2 = 0.4452295445267456
$ws.Cells.Item(4, 10).Value2 = 0.3751312191747254
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 21.059022
$ws.Cells.Item(4, 14).Value2 = 63.177066
$ws.Cells.Item(4, 15).Value2 = 0.521923552418633
$ws.Cells.Item(4, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(4, 17).Value2 = 808.4395759559941
$ws.Cells.Item(4, 18).Value2 = 4850.637455735964
$ws.Cells.Item(4, 19).Value2 = 0.232375785521129
$ws.Cells.Item(4, 20).Value2 = 0.1959126092713747

$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 7).Value2 = 38.38922700000001
$ws.Cells.Item(5, 8).Value2 = 76.77845400000001
$ws.Cells.Item(5, 9).Value2 = 0.4452295445267456
$ws.Cells.Item(5, 10).Value2 = 0.3751312191747254
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 13).Value2 = 0.131278
$ws.Cells.Item(5, 14).Value2 = 0.393834
$ws.Cells.Item(5, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(5, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(5, 17).Value2 = 5.039660942106001
$ws.Cells.Item(5, 18).Value2 = 30.23796565263601
$ws.Cells.Item(5, 19).Value2 = 0.001448587136270752
$ws.Cells.Item(5, 20).Value2 = 0.001221282522993116

$ws.Cells.Item(6, 5).Value2 = 2
$ws.Cells.Item(6, 7).Value2 = 38.38922700000001
$ws.Cells.Item(6, 8).Value2 = 76.77845400000001
$ws.Cells.Item(6, 9).Value2 = 0.4452295445267456
$ws.Cells.Item(6, 10).Value2 = 0.3751312191747254
$ws.Cells.Item(6, 11).Value2 = 2
$ws.Cells.Item(6, 13).Value2 = 0.0758675
$ws.Cells.Item(6, 14).Value2 = 0.151735
$ws.Cells.Item(6, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(6, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(6, 17).Value2 = 2.912494679422501
$ws.Cells.Item(6, 18).Value2 = 11.64997871769
$ws.Cells.Item(6, 19).Value2 = 0.0008371599549126381
$ws.Cells.Item(6, 20).Value2 = 0.0004705315021718806

$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 7).Value2 = 3.154770000000001
$ws.Cells.Item(7, 8).Value2 = 9.464310000000001
$ws.Cells.Item(7, 9).Value2 = 0.03658830666704076
$ws.Cells.Item(7, 10).Value2 = 0.04624159466596638
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 13).Value2 = 0.431063
$ws.Cells.Item(7, 14).Value2 = 1.293189
$ws.Cells.Item(7, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(7, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(7, 17).Value2 = 1.35990462051
$ws.Cells.Item(7, 18).Value2 = 12.23914158459
$ws.Cells.Item(7, 19).Value2 = 0.0003908874748630879
$ws.Cells.Item(7, 20).Value2 = 0.0004943272270829829

$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 7).Value2 = 3.154770000000001
$ws.Cells.Item(8, 8).Value2 = 9.464310000000001
$ws.Cells.Item(8, 9).Value2 = 0.03658830666704076
$ws.Cells.Item(8, 10).Value2 = 0.04624159466596638
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 13).Value2 = 18.65163266666667
$ws.Cells.Item(8, 14).Value2 = 55.954898
$ws.Cells.Item(8, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(8, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(8, 17).Value2 = 58.84161118782001
$ws.Cells.Item(8, 18).Value2 = 529.5745006903801
$ws.Cells.Item(8, 19).Value2 = 0.01691328087807865
$ws.Cells.Item(8, 20).Value2 = 0.02138900777075211

$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 7).Value2 = 3.154770000000001
$ws.Cells.Item(9, 8).Value2 = 9.464310000000001
$ws.Cells.Item(9, 9).Value2 = 0.03658830666704076
$ws.Cells.Item(9, 10).Value2 = 0.04624159466596638
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 13).Value2 = 21.059022
$ws.Cells.Item(9, 14).Value2 = 63.177066
$ws.Cells.Item(9, 15).Value2 = 0.521923552418633
$ws.Cells.Item(9, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(9, 17).Value2 = 66.43637083494
$ws.Cells.Item(9, 18).Value2 = 597.9273375144601
$ws.Cells.Item(9, 19).Value2 = 0.01909629899264427
$ws.Cells.Item(9, 20).Value2 = 0.02414971349974258

$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 7).Value2 = 3.154770000000001
$ws.Cells.Item(10, 8).Value2 = 9.464310000000001
$ws.Cells.Item(10, 9).Value2 = 0.03658830666704076
$ws.Cells.Item(10, 10).Value2 = 0.04624159466596638
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 13).Value2 = 0.131278
$ws.Cells.Item(10, 14).Value2 = 0.393834
$ws.Cells.Item(10, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(10, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(10, 17).Value2 = 0.4141518960600001
$ws.Cells.Item(10, 18).Value2 = 3.727367064540001
$ws.Cells.Item(10, 19).Value2 = 0.0001190427522776867
$ws.Cells.Item(10, 20).Value2 = 0.0001505447920999943

$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 7).Value2 = 3.154770000000001
$ws.Cells.Item(11, 8).Value2 = 9.464310000000001
$ws.Cells.Item(11, 9).Value2 = 0.03658830666704076
$ws.Cells.Item(11, 10).Value2 = 0.04624159466596638
$ws.Cells.Item(11, 11).Value2 = 2
$ws.Cells.Item(11, 13).Value2 = 0.0758675
$ws.Cells.Item(11, 14).Value2 = 0.151735
$ws.Cells.Item(11, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(11, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(11, 17).Value2 = 0.2393445129750001
$ws.Cells.Item(11, 18).Value2 = 1.43606707785
$ws.Cells.Item(11, 19).Value2 = 0.00006879656917707
$ws.Cells.Item(11, 20).Value2 = 0.00005800137628872224

$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 7).Value2 = 2.963441333333333
$ws.Cells.Item(12, 8).Value2 = 8.890324
$ws.Cells.Item(12, 9).Value2 = 0.03436932020203823
$ws.Cells.Item(12, 10).Value2 = 0.04343716117256439
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 13).Value2 = 0.431063
$ws.Cells.Item(12, 14).Value2 = 1.293189
$ws.Cells.Item(12, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(12, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(12, 17).Value2 = 1.277429911470667
$ws.Cells.Item(12, 18).Value2 = 11.496869203236
$ws.Cells.Item(12, 19).Value2 = 0.0003671811573241691
$ws.Cells.Item(12, 20).Value2 = 0.0004643475552670287

$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 7).Value2 = 2.963441333333333
$ws.Cells.Item(13, 8).Value2 = 8.890324
$ws.Cells.Item(13, 9).Value2 = 0.03436932020203823
$ws.Cells.Item(13, 10).Value2 = 0.04343716117256439
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 13).Value2 = 18.65163266666667
$ws.Cells.Item(13, 14).Value2 = 55.954898
$ws.Cells.Item(13, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(13, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(13, 17).Value2 = 55.27301917855023
$ws.Cells.Item(13, 18).Value2 = 497.457172606952
$ws.Cells.Item(13, 19).Value2 = 0.01588753399974469
$ws.Cells.Item(13, 20).Value2 = 0.02009181959598787

$ws.Cells.Item(14, 5).Value2 = 3
$ws.Cells.Item(14, 7).Value2 = 2.963441333333333
$ws.Cells.Item(14, 8).Value2 = 8.890324
$ws.Cells.Item(14, 9).Value2 = 0.03436932020203823
$ws.Cells.Item(14, 10).Value2 = 0.04343716117256439
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 13).Value2 = 21.059022
$ws.Cells.Item(14, 14).Value2 = 63.177066
$ws.Cells.Item(14, 15).Value2 = 0.521923552418633
$ws.Cells.Item(14, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(14, 17).Value2 = 62.407176234376
$ws.Cells.Item(14, 18).Value2 = 561.664586109384
$ws.Cells.Item(14, 19).Value2 = 0.01793815769406128
$ws.Cells.Item(14, 20).Value2 = 0.02268509564034625

$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 7).Value2 = 2.963441333333333
$ws.Cells.Item(15, 8).Value2 = 8.890324
$ws.Cells.Item(15, 9).Value2 = 0.03436932020203823
$ws.Cells.Item(15, 10).Value2 = 0.04343716117256439
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 13).Value2 = 0.131278
$ws.Cells.Item(15, 14).Value2 = 0.393834
$ws.Cells.Item(15, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(15, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(15, 17).Value2 = 0.3890346513573333
$ws.Cells.Item(15, 18).Value2 = 3.501311862216
$ws.Cells.Item(15, 19).Value2 = 0.0001118231162758165
$ws.Cells.Item(15, 20).Value2 = 0.0001414146386035104

$ws.Cells.Item(16, 5).Value2 = 3
$ws.Cells.Item(16, 7).Value2 = 2.963441333333333
$ws.Cells.Item(16, 8).Value2 = 8.890324
$ws.Cells.Item(16, 9).Value2 = 0.03436932020203823
$ws.Cells.Item(16, 10).Value2 = 0.04343716117256439
$ws.Cells.Item(16, 11).Value2 = 2
$ws.Cells.Item(16, 13).Value2 = 0.0758675
$ws.Cells.Item(16, 14).Value2 = 0.151735
$ws.Cells.Item(16, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(16, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(16, 17).Value2 = 0.2248288853566667
$ws.Cells.Item(16, 18).Value2 = 1.34897331214
$ws.Cells.Item(16, 19).Value2 = 0.00006462423463227278
$ws.Cells.Item(16, 20).Value2 = 0.00005448374235973443

$ws.Cells.Item(17, 5).Value2 = 3
$ws.Cells.Item(17, 7).Value2 = 7.153525333333334
$ws.Cells.Item(17, 8).Value2 = 21.460576
$ws.Cells.Item(17, 9).Value2 = 0.08296496373632466
$ws.Cells.Item(17, 10).Value2 = 0.1048540524021472
$ws.Cells.Item(17, 11).Value2 = 3
$ws.Cells.Item(17, 13).Value2 = 0.431063
$ws.Cells.Item(17, 14).Value2 = 1.293189
$ws.Cells.Item(17, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(17, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(17, 17).Value2 = 3.083620090762667
$ws.Cells.Item(17, 18).Value2 = 27.752580816864
$ws.Cells.Item(17, 19).Value2 = 0.0008863478015563085
$ws.Cells.Item(17, 20).Value2 = 0.001120900205686797

$ws.Cells.Item(18, 5).Value2 = 3
$ws.Cells.Item(18, 7).Value2 = 7.153525333333334
$ws.Cells.Item(18, 8).Value2 = 21.460576
$ws.Cells.Item(18, 9).Value2 = 0.08296496373632466
$ws.Cells.Item(18, 10).Value2 = 0.1048540524021472
$ws.Cells.Item(18, 11).Value2 = 3
$ws.Cells.Item(18, 13).Value2 = 18.65163266666667
$ws.Cells.Item(18, 14).Value2 = 55.954898
$ws.Cells.Item(18, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(18, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(18, 17).Value2 = 133.4249267890276
$ws.Cells.Item(18, 18).Value2 = 1200.824341101248
$ws.Cells.Item(18, 19).Value2 = 0.0383513166510135
$ws.Cells.Item(18, 20).Value2 = 0.04850014706078058

$ws.Cells.Item(19, 5).Value2 = 3
$ws.Cells.Item(19, 7).Value2 = 7.153525333333334
$ws.Cells.Item(19, 8).Value2 = 21.460576
$ws.Cells.Item(19, 9).Value2 = 0.08296496373632466
$ws.Cells.Item(19, 10).Value2 = 0.1048540524021472
$ws.Cells.Item(19, 11).Value2 = 3
$ws.Cells.Item(19, 13).Value2 = 21.059022
$ws.Cells.Item(19, 14).Value2 = 63.177066
$ws.Cells.Item(19, 15).Value2 = 0.521923552418633
$ws.Cells.Item(19, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(19, 17).Value2 = 150.646247372224
$ws.Cells.Item(19, 18).Value2 = 1355.816226350016
$ws.Cells.Item(19, 19).Value2 = 0.04330136859954563
$ws.Cells.Item(19, 20).Value2 = 0.05476012112234822

$ws.Cells.Item(20, 5).Value2 = 3
$ws.Cells.Item(20, 7).Value2 = 7.153525333333334
$ws.Cells.Item(20, 8).Value2 = 21.460576
$ws.Cells.Item(20, 9).Value2 = 0.08296496373632466
$ws.Cells.Item(20, 10).Value2 = 0.1048540524021472
$ws.Cells.Item(20, 11).Value2 = 3
$ws.Cells.Item(20, 13).Value2 = 0.131278
$ws.Cells.Item(20, 14).Value2 = 0.393834
$ws.Cells.Item(20, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(20, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(20, 17).Value2 = 0.9391004987093334
$ws.Cells.Item(20, 18).Value2 = 8.451904488384
$ws.Cells.Item(20, 19).Value2 = 0.0002699326239846822
$ws.Cells.Item(20, 20).Value2 = 0.0003413643416441481

$ws.Cells.Item(21, 5).Value2 = 3
$ws.Cells.Item(21, 7).Value2 = 7.153525333333334
$ws.Cells.Item(21, 8).Value2 = 21.460576
$ws.Cells.Item(21, 9).Value2 = 0.08296496373632466
$ws.Cells.Item(21, 10).Value2 = 0.1048540524021472
$ws.Cells.Item(21, 11).Value2 = 2
$ws.Cells.Item(21, 13).Value2 = 0.0758675
$ws.Cells.Item(21, 14).Value2 = 0.151735
$ws.Cells.Item(21, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(21, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(21, 17).Value2 = 0.5427200832266668
$ws.Cells.Item(21, 18).Value2 = 3.25632049936
$ws.Cells.Item(21, 19).Value2 = 0.0001559980602245455
$ws.Cells.Item(21, 20).Value2 = 0.0001315196716874998

$ws.Cells.Item(22, 5).Value2 = 3
$ws.Cells.Item(22, 7).Value2 = 18.952291
$ws.Cells.Item(22, 8).Value2 = 56.85687299999999
$ws.Cells.Item(22, 9).Value2 = 0.2198043708894773
$ws.Cells.Item(22, 10).Value2 = 0.2777965298305241
$ws.Cells.Item(22, 11).Value2 = 3
$ws.Cells.Item(22, 13).Value2 = 0.431063
$ws.Cells.Item(22, 14).Value2 = 1.293189
$ws.Cells.Item(22, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(22, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(22, 17).Value2 = 8.169631415332999
$ws.Cells.Item(22, 18).Value2 = 73.52668273799699
$ws.Cells.Item(22, 19).Value2 = 0.002348257772154682
$ws.Cells.Item(22, 20).Value2 = 0.00296967241887674

$ws.Cells.Item(23, 5).Value2 = 3
$ws.Cells.Item(23, 7).Value2 = 18.952291
$ws.Cells.Item(23, 8).Value2 = 56.85687299999999
$ws.Cells.Item(23, 9).Value2 = 0.2198043708894773
$ws.Cells.Item(23, 10).Value2 = 0.2777965298305241
$ws.Cells.Item(23, 11).Value2 = 3
$ws.Cells.Item(23, 13).Value2 = 18.65163266666667
$ws.Cells.Item(23, 14).Value2 = 55.954898
$ws.Cells.Item(23, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(23, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(23, 17).Value2 = 353.4911699237726
$ws.Cells.Item(23, 18).Value2 = 3181.420529313953
$ws.Cells.Item(23, 19).Value2 = 0.1016065896930939
$ws.Cells.Item(23, 20).Value2 = 0.1284945335072145

$ws.Cells.Item(24, 5).Value2 = 3
$ws.Cells.Item(24, 7).Value2 = 18.952291
$ws.Cells.Item(24, 8).Value2 = 56.85687299999999
$ws.Cells.Item(24, 9).Value2 = 0.2198043708894773
$ws.Cells.Item(24, 10).Value2 = 0.2777965298305241
$ws.Cells.Item(24, 11).Value2 = 3
$ws.Cells.Item(24, 13).Value2 = 21.059022
$ws.Cells.Item(24, 14).Value2 = 63.177066
$ws.Cells.Item(24, 15).Value2 = 0.521923552418633
$ws.Cells.Item(24, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(24, 17).Value2 = 399.1167131194019
$ws.Cells.Item(24, 18).Value2 = 3592.050418074617
$ws.Cells.Item(24, 19).Value2 = 0.1147210780917788
$ws.Cells.Item(24, 20).Value2 = 0.1450794821219137

$ws.Cells.Item(25, 5).Value2 = 3
$ws.Cells.Item(25, 7).Value2 = 18.952291
$ws.Cells.Item(25, 8).Value2 = 56.85687299999999
$ws.Cells.Item(25, 9).Value2 = 0.2198043708894773
$ws.Cells.Item(25, 10).Value2 = 0.2777965298305241
$ws.Cells.Item(25, 11).Value2 = 3
$ws.Cells.Item(25, 13).Value2 = 0.131278
$ws.Cells.Item(25, 14).Value2 = 0.393834
$ws.Cells.Item(25, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(25, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(25, 17).Value2 = 2.488018857898
$ws.Cells.Item(25, 18).Value2 = 22.392169721082
$ws.Cells.Item(25, 19).Value2 = 0.0007151497201404953
$ws.Cells.Item(25, 20).Value2 = 0.000904398326475018

$ws.Cells.Item(26, 5).Value2 = 3
$ws.Cells.Item(26, 7).Value2 = 18.952291
$ws.Cells.Item(26, 8).Value2 = 56.85687299999999
$ws.Cells.Item(26, 9).Value2 = 0.2198043708894773
$ws.Cells.Item(26, 10).Value2 = 0.2777965298305241
$ws.Cells.Item(26, 11).Value2 = 2
$ws.Cells.Item(26, 13).Value2 = 0.0758675
$ws.Cells.Item(26, 14).Value2 = 0.151735
$ws.Cells.Item(26, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(26, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(26, 17).Value2 = 1.4378629374425
$ws.Cells.Item(26, 18).Value2 = 8.627177624655
$ws.Cells.Item(26, 19).Value2 = 0.0004132956123094428
$ws.Cells.Item(26, 20).Value2 = 0.0003484434560441376

$ws.Cells.Item(27, 5).Value2 = 2
$ws.Cells.Item(27, 7).Value2 = 15.6101945
$ws.Cells.Item(27, 8).Value2 = 31.220389
$ws.Cells.Item(27, 9).Value2 = 0.1810434939783733
$ws.Cells.Item(27, 10).Value2 = 0.1525394427540724
$ws.Cells.Item(27, 11).Value2 = 3
$ws.Cells.Item(27, 13).Value2 = 0.431063
$ws.Cells.Item(27, 14).Value2 = 1.293189
$ws.Cells.Item(27, 15).Value2 = 0.01068339889080477
$ws.Cells.Item(27, 16).Value2 = 0.01069009904727195
$ws.Cells.Item(27, 17).Value2 = 6.728977271753499
$ws.Cells.Item(27, 18).Value2 = 40.37386363052099
$ws.Cells.Item(27, 19).Value2 = 0.001934159862755973
$ws.Cells.Item(27, 20).Value2 = 0.001630661751656704

$ws.Cells.Item(28, 5).Value2 = 2
$ws.Cells.Item(28, 7).Value2 = 15.6101945
$ws.Cells.Item(28, 8).Value2 = 31.220389
$ws.Cells.Item(28, 9).Value2 = 0.1810434939783733
$ws.Cells.Item(28, 10).Value2 = 0.1525394427540724
$ws.Cells.Item(28, 11).Value2 = 3
$ws.Cells.Item(28, 13).Value2 = 18.65163266666667
$ws.Cells.Item(28, 14).Value2 = 55.954898
$ws.Cells.Item(28, 15).Value2 = 0.4622591865754301
$ws.Cells.Item(28, 16).Value2 = 0.4625490951438647
$ws.Cells.Item(28, 17).Value2 = 291.1556136692203
$ws.Cells.Item(28, 18).Value2 = 1746.933682015322
$ws.Cells.Item(28, 19).Value2 = 0.08368901826121664
$ws.Cells.Item(28, 20).Value2 = 0.07055698121964556

$ws.Cells.Item(29, 5).Value2 = 2
$ws.Cells.Item(29, 7).Value2 = 15.6101945
$ws.Cells.Item(29, 8).Value2 = 31.220389
$ws.Cells.Item(29, 9).Value2 = 0.1810434939783733
$ws.Cells.Item(29, 10).Value2 = 0.1525394427540724
$ws.Cells.Item(29, 11).Value2 = 3
$ws.Cells.Item(29, 13).Value2 = 21.059022
$ws.Cells.Item(29, 14).Value2 = 63.177066
$ws.Cells.Item(29, 15).Value2 = 0.521923552418633
$ws.Cells.Item(29, 16).Value2 = 0.5222508798451249
$ws.Cells.Item(29, 17).Value2 = 328.735429399779
$ws.Cells.Item(29, 18).Value2 = 1972.412576398674
$ws.Cells.Item(29, 19).Value2 = 0.094490863519474
$ws.Cells.Item(29, 20).Value2 = 0.07966385818939939

$ws.Cells.Item(30, 5).Value2 = 2
$ws.Cells.Item(30, 7).Value2 = 15.6101945
$ws.Cells.Item(30, 8).Value2 = 31.220389
$ws.Cells.Item(30, 9).Value2 = 0.1810434939783733
$ws.Cells.Item(30, 10).Value2 = 0.1525394427540724
$ws.Cells.Item(30, 11).Value2 = 3
$ws.Cells.Item(30, 13).Value2 = 0.131278
$ws.Cells.Item(30, 14).Value2 = 0.393834
$ws.Cells.Item(30, 15).Value2 = 0.003253573699406046
$ws.Cells.Item(30, 16).Value2 = 0.003255614197293127
$ws.Cells.Item(30, 17).Value2 = 2.049275113571
$ws.Cells.Item(30, 18).Value2 = 12.295650681426
$ws.Cells.Item(30, 19).Value2 = 0.0005890383504566123
$ws.Cells.Item(30, 20).Value2 = 0.0004966095754773404

$ws.Cells.Item(31, 5).Value2 = 2
$ws.Cells.Item(31, 7).Value2 = 15.6101945
$ws.Cells.Item(31, 8).Value2 = 31.220389
$ws.Cells.Item(31, 9).Value2 = 0.1810434939783733
$ws.Cells.Item(31, 10).Value2 = 0.1525394427540724
$ws.Cells.Item(31, 11).Value2 = 2
$ws.Cells.Item(31, 13).Value2 = 0.0758675
$ws.Cells.Item(31, 14).Value2 = 0.151735
$ws.Cells.Item(31, 15).Value2 = 0.001880288415726079
$ws.Cells.Item(31, 16).Value2 = 0.001254311766445438
$ws.Cells.Item(31, 17).Value2 = 1.18430643122875
$ws.Cells.Item(31, 18).Value2 = 4.737225724915
$ws.Cells.Item(31, 19).Value2 = 0.0003404139844701095
$ws.Cells.Item(31, 20).Value2 = 0.0001913320178934634
